$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Entry order reproduces the shared-strings table order seen in the
# target workbook.
$ws.Range("A15").Value = "F"

$ws.Range("C14").Value = "jpc mem"
$ws.Range("D14").Value = "If (overflow) goto mem"
$ws.Range("B14").Value = "Jump if overflow"
$ws.Range("A14").Value = "E0XY"

$ws.Range("A13").Value = "DRST"
$ws.Range("D13").Value = "reg = reg1 x reg2"
$ws.Range("C13").Value = "mulf reg,reg1,reg2"
$ws.Range("B13").Value = "Multiplicar fp"

$ws.Range("C16").Select()
